# Generate Report for Handback
#
# The handback-report generator has processed a new handback for the
# "605d0248-9444-4adb-9df3-615e76a5b920" source file (row 6 in both the
# zh-cn and de-de status tables). For each locale it fills in:
#   - Latest Target File   (column I) -> the source-file's handback .md name
#                                        (also turned into a hyperlink, like A6)
#   - Latest Handback File (column J) -> the locale's generated .xlf file name
#   - Latest Handback DateTime (column K) -> the timestamp the handback report ran
#   - Error Detail (column P) -> a warning that the handback file isn't the
#                                 tip-of-branch version
# It also widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$handbackMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3d7a7101474a7a2c3e2e96a46eb63b411f8c603/e2e/605d0248-9444-4adb-9df3-615e76a5b920.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/981dc4e309cc7c784f6b322ca99b1d7eb7c07d81/e2e/605d0248-9444-4adb-9df3-615e76a5b920.md."
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/981dc4e309cc7c784f6b322ca99b1d7eb7c07d81/e2e/605d0248-9444-4adb-9df3-615e76a5b920.md"
$targetFileDisplay = "605d0248-9444-4adb-9df3-615e76a5b920.md"

# -- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Column P ("Error Detail") is about to receive a long message - widen it.
$wsZh.Columns.Item(16).ColumnWidth = 39.1667

$wsZh.Range("I6").Value = $targetFileDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $targetFileUrl, $null, $null, $targetFileDisplay)
$wsZh.Range("I6").Font.Underline = $true
$wsZh.Range("I6").Font.Color = 15570276

$wsZh.Range("J6").Value = "605d0248-9444-4adb-9df3-615e76a5b920.39655294ad9aa24a25a76b554b4c6c229323a21f.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-09-04 16:47:50"
$wsZh.Range("P6").Value = $handbackMessage

# -- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.1667

$wsDe.Range("I6").Value = $targetFileDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $targetFileUrl, $null, $null, $targetFileDisplay)
$wsDe.Range("I6").Font.Underline = $true
$wsDe.Range("I6").Font.Color = 15570276

$wsDe.Range("J6").Value = "605d0248-9444-4adb-9df3-615e76a5b920.39655294ad9aa24a25a76b554b4c6c229323a21f.de-de.xlf"
$wsDe.Range("K6").Value = "2016-09-04 16:47:57"
$wsDe.Range("P6").Value = $handbackMessage
